# BL Audit Form - "last report 19-02-25" update
# Updates the report date, the payment note, several stock-quantity inputs,
# the two receivable/due figures, clears the old bank-guarantee figure and
# bumps the commission figure, then moves the on-screen selection to match
# where the user ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Report date (B1) and payment note (F34) ---------------------------
$ws.Range("B1").Value = "19.02.2025"
$ws.Range("F34").Value = "20.02.2025 payment "

# --- Stock table (rows 9-17): quantities received ----------------------
$ws.Range("C9").Value  = 322015
$ws.Range("C10").Value = 1140
$ws.Range("C11").Value = 290
$ws.Range("C12").Value = 5610
$ws.Range("C14").Value = 26
$ws.Range("C16").Value = 53
$ws.Range("C17").Value = 77

# --- Receivables / dues --------------------------------------------------
$ws.Range("E22").Value = 50089
$ws.Range("E23").Value = 55434

# --- Bank guarantee figure cleared ---------------------------------------
$ws.Range("E31").ClearContents()

# --- Commission figure ----------------------------------------------------
$ws.Range("E34").Value = 120000

# --- Recalculate so dependent formulas (E9:E25, E32, E39, E40) refresh ---
$excel.Calculate()

# --- Scroll / selection position, matching where the editor ended up ----
$ws.Range("E35").Select()
$excel.ActiveWindow.ScrollRow = 23
